$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column, copying the header formatting from an existing
# header cell (G1) so H1 reuses the same bold/border/alignment style.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-7
$saveValues = @(1, 0, 0, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
